$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.928768402600219
$ws.Range("D2").Value = 5.649667461805564
$ws.Range("E2").Value = 16.72339935239674
$ws.Range("F2").Value = 38.60592147596483
$ws.Range("G2").Value = 59.47479272158546
$ws.Range("H2").Value = 18.02939739762821

$ws.Range("C3").Value = 4.750715973229781
$ws.Range("D3").Value = 5.449517218026375
$ws.Range("E3").Value = 15.74130727870119
$ws.Range("F3").Value = 36.83191499190448
$ws.Range("G3").Value = 56.11544429333878
$ws.Range("H3").Value = 17.50004222806566

$ws.Range("C4").Value = 4.640201828005192
$ws.Range("D4").Value = 5.325647932214536
$ws.Range("E4").Value = 15.11449299435107
$ws.Range("F4").Value = 35.71611540401494
$ws.Range("G4").Value = 53.97312806182247
$ws.Range("H4").Value = 17.17499675511617

$ws.Range("C5").Value = 4.594943887929193
$ws.Range("D5").Value = 5.275008664708693
$ws.Range("E5").Value = 14.85333479962914
$ws.Range("F5").Value = 35.25546152999916
$ws.Range("G5").Value = 53.08108181568053
$ws.Range("H5").Value = 17.04274510246151

$ws.Range("C6").Value = 4.587417753661993
$ws.Range("D6").Value = 5.266592867242308
$ws.Range("E6").Value = 14.80963250433649
$ws.Range("F6").Value = 35.17863344159592
$ws.Range("G6").Value = 52.93184151972932
$ws.Range("H6").Value = 17.02080368120299

$ws.Range("C7").Value = 4.639592259107098
$ws.Range("D7").Value = 5.324965532006909
$ws.Range("E7").Value = 15.11099372557861
$ws.Range("F7").Value = 35.70992596150554
$ws.Range("G7").Value = 53.96117322644588
$ws.Range("H7").Value = 17.17321203841497

$ws.Range("C8").Value = 4.867666770924492
$ws.Range("D8").Value = 5.580906729155222
$ws.Range("E8").Value = 16.38985453637523
$ws.Range("F8").Value = 38.00015848500796
$ws.Range("G8").Value = 58.33356360182393
$ws.Range("H8").Value = 17.84700391936297

$ws.Range("C9").Value = 5.302440707935573
$ws.Range("D9").Value = 6.071699395578595
$ws.Range("E9").Value = 18.86970573823709
$ws.Range("F9").Value = 42.25359736752895
$ws.Range("G9").Value = 66.24093425936105
$ws.Range("H9").Value = 19.15937317617588

$ws.Range("C10").Value = 5.610646901156922
$ws.Range("D10").Value = 6.421499068240283
$ws.Range("E10").Value = 20.6068665749502
$ws.Range("F10").Value = 45.20348272566346
$ws.Range("G10").Value = 71.60792151330536
$ws.Range("H10").Value = 20.10724551707758

$ws.Range("C11").Value = 5.747787228625389
$ws.Range("D11").Value = 6.577581296935055
$ws.Range("E11").Value = 21.35638126211074
$ws.Range("F11").Value = 46.50257961559647
$ws.Range("G11").Value = 73.9482740044039
$ws.Range("H11").Value = 20.53288917723165

$ws.Range("C12").Value = 5.79923357677165
$ws.Range("D12").Value = 6.636197990370094
$ws.Range("E12").Value = 21.63440586590772
$ws.Range("F12").Value = 46.98803801519092
$ws.Range("G12").Value = 74.81966114465952
$ws.Range("H12").Value = 20.6931227634408

$ws.Range("C13").Value = 5.788175943633451
$ws.Range("D13").Value = 6.623596289287762
$ws.Range("E13").Value = 21.57478535990522
$ws.Range("F13").Value = 46.88377867963652
$ws.Range("G13").Value = 74.63265710415702
$ws.Range("H13").Value = 20.65865801399359

$ws.Range("C14").Value = 5.752029722615116
$ws.Range("D14").Value = 6.58241377987238
$ws.Range("E14").Value = 21.37937043541524
$ws.Range("F14").Value = 46.54265012011395
$ws.Range("G14").Value = 74.02026275025013
$ws.Range("H14").Value = 20.54609150817867

$ws.Range("C15").Value = 5.729824625547359
$ws.Range("D15").Value = 6.557123326602148
$ws.Range("E15").Value = 21.25891946079912
$ws.Range("F15").Value = 46.33284632113646
$ws.Range("G15").Value = 73.64321121104929
$ws.Range("H15").Value = 20.47701354564039

$ws.Range("C16").Value = 5.601618809704126
$ws.Range("D16").Value = 6.411233009933231
$ws.Range("E16").Value = 20.55706926943454
$ws.Range("F16").Value = 45.11769302263878
$ws.Range("G16").Value = 71.45291064196773
$ws.Range("H16").Value = 20.07930384196964

$ws.Range("C17").Value = 5.522149913736917
$ws.Range("D17").Value = 6.320916434563105
$ws.Range("E17").Value = 20.11611526749959
$ws.Range("F17").Value = 44.36101940554379
$ws.Range("G17").Value = 70.08308494553567
$ws.Range("H17").Value = 19.83378887045362

$ws.Range("C18").Value = 5.476155396883597
$ws.Range("D18").Value = 6.268684870150421
$ws.Range("E18").Value = 19.85865322088635
$ws.Range("F18").Value = 43.92178569123348
$ws.Range("G18").Value = 69.28569292124584
$ws.Range("H18").Value = 19.69206097939126

$ws.Range("C19").Value = 5.460534774780904
$ws.Range("D19").Value = 6.25095305674084
$ws.Range("E19").Value = 19.77082028624357
$ws.Range("F19").Value = 43.7723899325955
$ws.Range("G19").Value = 69.01408814294487
$ws.Range("H19").Value = 19.64399087573018

$ws.Range("C20").Value = 5.530639469913202
$ws.Range("D20").Value = 6.330560553955038
$ws.Range("E20").Value = 20.16345227902731
$ws.Range("F20").Value = 44.44198667598485
$ws.Range("G20").Value = 70.2298914315606
$ws.Range("H20").Value = 19.85997871186778

$ws.Range("C21").Value = 5.762660256602844
$ws.Range("D21").Value = 6.594523703688522
$ws.Range("E21").Value = 21.4369255030305
$ws.Range("F21").Value = 46.64302604969184
$ws.Range("G21").Value = 74.20054304427929
$ws.Range("H21").Value = 20.57918188555966

$ws.Range("C22").Value = 5.911451246215811
$ws.Range("D22").Value = 6.764173959020208
$ws.Range("E22").Value = 22.23543168887845
$ws.Range("F22").Value = 48.04364248591192
$ws.Range("G22").Value = 76.70892311719525
$ws.Range("H22").Value = 21.04363524480705

$ws.Range("C23").Value = 5.83231269242143
$ws.Range("D23").Value = 6.67390562883974
$ws.Range("E23").Value = 21.81232618168076
$ws.Range("F23").Value = 47.29966667446963
$ws.Range("G23").Value = 75.3781675039785
$ws.Range("H23").Value = 20.79630489046479

$ws.Range("C24").Value = 5.526802293663678
$ws.Range("D24").Value = 6.326201401030127
$ws.Range("E24").Value = 20.14206352909109
$ws.Range("F24").Value = 44.40539446720129
$ws.Range("G24").Value = 70.16355090943041
$ws.Range("H24").Value = 19.84814007350923

$ws.Range("C25").Value = 5.186557253758299
$ws.Range("D25").Value = 5.940557437116683
$ws.Range("E25").Value = 18.19349308199893
$ws.Range("F25").Value = 41.13180146318761
$ws.Range("G25").Value = 64.17740294065266
$ws.Range("H25").Value = 18.80647406668156
